$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert blank rows at both insertion points first (working bottom-to-top so
# row numbers stay valid), then fill in the values in "A-before-NU" order so
# new shared-string entries are created in that order.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Rows 5 & 6: new "Manufacturer_AT" / "AuthorisedRep_AT" entries, inserted
# right before the "Noor.Uddin.Business" row.
$ws.Cells.Item(5, 1).Value = "Manufacturer_AT"
$ws.Cells.Item(5, 2).Value = "MHRA12345A"
$ws.Cells.Item(5, 3).Value = "yes"
$ws.Cells.Item(5, 4).Value = "AT"

$ws.Cells.Item(6, 1).Value = "AuthorisedRep_AT"
$ws.Cells.Item(6, 2).Value = "MHRA12345A"
$ws.Cells.Item(6, 3).Value = "yes"
$ws.Cells.Item(6, 4).Value = "AT"

# Rows 10 & 11: new "Manufacturer_NU" / "AuthorisedRep_NU" entries, inserted
# right after the "Noor.Uddin.AuthorisedRep" row.
$ws.Cells.Item(10, 1).Value = "Manufacturer_NU"
$ws.Cells.Item(10, 2).Value = "MHRA12345A"
$ws.Cells.Item(10, 3).Value = "yes"
$ws.Cells.Item(10, 4).Value = "AT"

$ws.Cells.Item(11, 1).Value = "AuthorisedRep_NU"
$ws.Cells.Item(11, 2).Value = "MHRA12345A"
$ws.Cells.Item(11, 3).Value = "yes"
$ws.Cells.Item(11, 4).Value = "AT"

# Selection moves to A11 in the edited workbook.
$ws.Range("A11").Select()
